# Apply the commit's changes:
#  1. Rename the shared-string test data "bala" -> "saul1" and
#     "bala1243" -> "saul123" on the InvalidLoginTest sheet (A2/B2).
#  2. Update the saved sheet selection on that sheet from a single cell
#     (B2) to a full-row style selection (A4:XFD1048576, active cell A4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLoginTest")
$ws.Activate()

$ws.Range("A2").Value = "saul1"
$ws.Range("B2").Value = "saul123"

$ws.Range("A4:XFD1048576").Select()
